$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Question Source"
$ws.Range("D1").Value = "Formula"

# Data rows
$ws.Range("B2").Value = "Logarithms"

$ws.Range("C2").Value = "200604003003"
$ws.Range("D2").Value = "\log_{a}b"

$ws.Range("C3").Value = "200604003003"
$ws.Range("D3").Value = "\frac{1}{\log_{a}b}"

$ws.Range("C4").Value = "200604003003"
$ws.Range("D4").Value = "2\log _{9} x + 1"

$ws.Range("C5").Value = "200604003003"
$ws.Range("D5").Value = "2\log _{x} 3"

$ws.Range("C6").Value = "20040400105"
$ws.Range("D6").Value = "\log_{16}(3x-1)"

$ws.Range("C7").Value = "20040400105"
$ws.Range("D7").Value = "\log_{4}(3x)+\log_{4}(0.5)"

$ws.Range("C8").Value = "20030400103"
$ws.Range("D8").Value = "\log_{2}x-\log_{4}(x-4)"

$ws.Range("C9").Value = "20080400104"
$ws.Range("D9").Value = "2 + \log_{3}(3x-7)"

$ws.Range("C10").Value = "20080400104"
$ws.Range("D10").Value = "\log_{3}(2x-3)"

$ws.Range("C11").Value = "19980200105"
$ws.Range("D11").Value = "\log_{3}(x-1)=2"

$ws.Range("C12").Value = "20010400107"
$ws.Range("D12").Value = "\log_{4}y +\log_{2}y"

$ws.Range("C13").Value = "20020200110"
$ws.Range("D13").Value = "\log_{4}y +\log_{2}y"

# apply numbers-as-text format (style index 1) to C2:C13, matching original C2:C3
$ws.Range("C2:C13").NumberFormat = "@"

# selection
$ws.Range("D14").Select()
